$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 1911.8572
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 1730.5
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 5191.5
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -5731.5
# Row 73
$ws.Range("H73").Value = 1911.8572
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 1730.5
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 5191.5
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -7063.5
# Row 132
$ws.Range("H132").Value = 9017642
$ws.Range("I132").Value = 15158775
$ws.Range("J132").Value = 10647
$ws.Range("K132").Value = 45476325
$ws.Range("L132").Value = 31941
$ws.Range("M132").Value = -45473795
$ws.Range("N132").Value = -37001
# Row 138
$ws.Range("H138").Value = 1477.0476
$ws.Range("J138").Value = 1668.2821
$ws.Range("L138").Value = 5004.846299999999
$ws.Range("N138").Value = -15284.8463

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5328.283
$ws.Range("I32").Value = 4772.255
$ws.Range("K32").Value = 4772.255
$ws.Range("M32").Value = -4485.255
# Row 45
$ws.Range("H45").Value = 1183.7693
$ws.Range("I45").Value = 1021
$ws.Range("K45").Value = 1021
$ws.Range("M45").Value = -644
# Row 61
$ws.Range("H61").Value = 83335900
$ws.Range("I61").Value = 90911624
$ws.Range("K61").Value = 90911624
$ws.Range("M61").Value = -90911412
# Row 122
$ws.Range("H122").Value = 956.8
$ws.Range("I122").Value = 986.25
$ws.Range("J122").Value = 839
$ws.Range("K122").Value = 2958.75
$ws.Range("L122").Value = 2517
$ws.Range("M122").Value = -508.75
$ws.Range("N122").Value = -7417
# Row 132
$ws.Range("H132").Value = 2554.2173
$ws.Range("I132").Value = 2042.4706
$ws.Range("J132").Value = 4004.1667
$ws.Range("K132").Value = 6127.4118
$ws.Range("L132").Value = 12012.5001
$ws.Range("M132").Value = -3597.4118
$ws.Range("N132").Value = -17072.5001
# Row 136
$ws.Range("H136").Value = 83335900
$ws.Range("I136").Value = 90911624
$ws.Range("K136").Value = 272734872
$ws.Range("M136").Value = -272732322

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 111111900
$ws.Range("I99").Value = 166667120
$ws.Range("J99").Value = 1466.6666
$ws.Range("K99").Value = 166667120
$ws.Range("L99").Value = 1466.6666
$ws.Range("M99").Value = -166665622
$ws.Range("N99").Value = -4462.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 62501296
$ws.Range("I16").Value = 90910370
$ws.Range("J16").Value = 1334
$ws.Range("K16").Value = 90910370
$ws.Range("L16").Value = 1334
$ws.Range("M16").Value = -90910083
$ws.Range("N16").Value = -1908
# Row 31
$ws.Range("H31").Value = 1206.8438
$ws.Range("I31").Value = 1147.4073
$ws.Range("J31").Value = 1527.8
$ws.Range("K31").Value = 1147.4073
$ws.Range("L31").Value = 1527.8
$ws.Range("M31").Value = -852.4073000000001
$ws.Range("N31").Value = -2117.8
# Row 34
$ws.Range("H34").Value = 1206.8438
$ws.Range("I34").Value = 1147.4073
$ws.Range("J34").Value = 1527.8
$ws.Range("K34").Value = 1147.4073
$ws.Range("L34").Value = 1527.8
$ws.Range("M34").Value = -945.4073000000001
$ws.Range("N34").Value = -1931.8
# Row 113
$ws.Range("H113").Value = 62501296
$ws.Range("I113").Value = 90910370
$ws.Range("J113").Value = 1334
$ws.Range("K113").Value = 90910370
$ws.Range("L113").Value = 1334
$ws.Range("M113").Value = -90908200
$ws.Range("N113").Value = -5674
# Row 132
$ws.Range("H132").Value = 10055.538
$ws.Range("I132").Value = 18204
$ws.Range("J132").Value = 3071.1428
$ws.Range("K132").Value = 54612
$ws.Range("L132").Value = 9213.428400000001
$ws.Range("M132").Value = -52082
$ws.Range("N132").Value = -14273.4284

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 3293.2
$ws.Range("J39").Value = 3257
$ws.Range("L39").Value = 9771
$ws.Range("N39").Value = -10359
# Row 55
$ws.Range("H55").Value = 2469.2307
$ws.Range("J55").Value = 3060
$ws.Range("L55").Value = 9180
$ws.Range("N55").Value = -9534
# Row 107
$ws.Range("H107").Value = 14842.714
$ws.Range("I107").Value = 660
$ws.Range("J107").Value = 20515.8
$ws.Range("K107").Value = 1980
$ws.Range("L107").Value = 61547.39999999999
$ws.Range("M107").Value = -60
$ws.Range("N107").Value = -65387.39999999999
# Row 131
$ws.Range("H131").Value = 16131922
$ws.Range("I131").Value = 142857490
$ws.Range("J131").Value = 3212.9272
$ws.Range("K131").Value = 428572470
$ws.Range("L131").Value = 9638.7816
$ws.Range("M131").Value = -428567430
$ws.Range("N131").Value = -19718.7816
# Row 140
$ws.Range("H140").Value = 21465.623
$ws.Range("I140").Value = 61009.883
$ws.Range("J140").Value = 2791.9443
$ws.Range("K140").Value = 183029.649
$ws.Range("L140").Value = 8375.832900000001
$ws.Range("M140").Value = -177849.649
$ws.Range("N140").Value = -18735.8329

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4558.4165
$ws.Range("I132").Value = 4633.6665
$ws.Range("K132").Value = 13900.9995
$ws.Range("M132").Value = -11370.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 746.3077
$ws.Range("I16").Value = 746.3077
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 746.3077
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -576.3077
$ws.Range("N16").ClearContents()
# Row 40
$ws.Range("H40").Value = 3624.5625
$ws.Range("I40").Value = 2318.182
$ws.Range("J40").Value = 6498.6
$ws.Range("K40").Value = 2318.182
$ws.Range("L40").Value = 6498.6
$ws.Range("M40").Value = -2182.182
$ws.Range("N40").Value = -6770.6
# Row 61
$ws.Range("H61").Value = 1087
$ws.Range("I61").Value = 1025
$ws.Range("J61").Value = 1335
$ws.Range("K61").Value = 1025
$ws.Range("L61").Value = 1335
$ws.Range("M61").Value = -823
$ws.Range("N61").Value = -1739
# Row 113
$ws.Range("H113").Value = 1087
$ws.Range("I113").Value = 1025
$ws.Range("J113").Value = 1335
$ws.Range("K113").Value = 1025
$ws.Range("L113").Value = 1335
$ws.Range("M113").Value = 1145
$ws.Range("N113").Value = -5675
# Row 122
$ws.Range("H122").Value = 27794512
$ws.Range("I122").Value = 41684350
$ws.Range("K122").Value = 125053050
$ws.Range("M122").Value = -125050600
# Row 128
$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 33340966
$ws.Range("I62").Value = 38467348
$ws.Range("J62").Value = 19501
$ws.Range("K62").Value = 38467348
$ws.Range("L62").Value = 19501
$ws.Range("M62").Value = -38466724
$ws.Range("N62").Value = -20749
# Row 65
$ws.Range("H65").Value = 33340966
$ws.Range("I65").Value = 38467348
$ws.Range("J65").Value = 19501
$ws.Range("K65").Value = 192336740
$ws.Range("L65").Value = 97505
$ws.Range("M65").Value = -192333620
$ws.Range("N65").Value = -103745
# Row 122
$ws.Range("H122").Value = 9617218
$ws.Range("I122").Value = 10418519
$ws.Range("K122").Value = 31255557
$ws.Range("M122").Value = -31253107
# Row 132
$ws.Range("H132").Value = 2557.125
$ws.Range("I132").Value = 1990.6666
$ws.Range("K132").Value = 5971.9998
$ws.Range("M132").Value = -3441.9998
# Row 136
$ws.Range("H136").Value = 909.4483
$ws.Range("I136").Value = 760.7619
$ws.Range("J136").Value = 1299.75
$ws.Range("K136").Value = 2282.2857
$ws.Range("L136").Value = 3899.25
$ws.Range("M136").Value = 267.7143000000001
$ws.Range("N136").Value = -8999.25
